# Schema update: multiple contributor roles, misc tidying
#
# - Contributor / Contributor1 sheets: rename the
#   "contributor__contributor_role" column header (F1) to
#   "contributor__contributor_roles".
# - Dataset / Dataset1 sheets: rename "dataset__title" (A1) to
#   "dataset__titles", and insert a new "dataset__submission_date"
#   column right before the existing "dataset__access_date" column
#   (shifting the remaining headers one column to the right).

$wb = $excel.ActiveWorkbook

foreach ($name in @("Contributor", "Contributor1")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F1").Value = "contributor__contributor_roles"
}

foreach ($name in @("Dataset", "Dataset1")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A1").Value = "dataset__titles"
    $ws.Range("C1").EntireColumn.Insert()
    $ws.Range("C1").Value = "dataset__submission_date"
}
